# Add the 3rd order/tier input row (School -> school -> tier 3) to Sheet1,
# mirroring the existing logs_type/type/current.tier rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A25").Value = "School"
$ws.Range("B25").Value = "school"
$ws.Range("C25").Value = 3

# Reflect the author's final view/selection state on the sheet (best effort;
# purely cosmetic window state).
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 8
$ws.Range("F26").Select() | Out-Null
